$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.926203513158953
$ws.Range("D2").Value = 2.696935350347433
$ws.Range("E2").Value = 9.551055133018757
$ws.Range("F2").Value = 20.21922706236545
$ws.Range("G2").Value = 24.03770515646365
$ws.Range("H2").Value = 11.68526316216086
$ws.Range("M2").Value = 19.2229271115908
$ws.Range("N2").Value = 17.18717348326094
$ws.Range("O2").Value = 17.36657399609888
$ws.Range("C3").Value = 4.750295668971268
$ws.Range("D3").Value = 2.691275868332867
$ws.Range("E3").Value = 9.730900890615883
$ws.Range("F3").Value = 19.82259608954563
$ws.Range("G3").Value = 23.15132912399772
$ws.Range("H3").Value = 11.64445299993882
$ws.Range("M3").Value = 18.43818122716748
$ws.Range("N3").Value = 16.88879049546276
$ws.Range("O3").Value = 17.1446210889951
$ws.Range("C4").Value = 4.640462652753223
$ws.Range("D4").Value = 2.687756574144641
$ws.Range("E4").Value = 9.846501334249391
$ws.Range("F4").Value = 19.58242330420264
$ws.Range("G4").Value = 22.59898279841709
$ws.Range("H4").Value = 11.62248522319309
$ws.Range("M4").Value = 17.93973538741299
$ws.Range("N4").Value = 16.7048604796733
$ws.Range("O4").Value = 17.01289557268333
$ws.Range("C5").Value = 4.595326335925164
$ws.Range("D5").Value = 2.686311698858806
$ws.Range("E5").Value = 9.894915050933371
$ws.Range("F5").Value = 19.48555034359608
$ws.Range("G5").Value = 22.37228249242559
$ws.Range("H5").Value = 11.61431598305299
$ws.Range("M5").Value = 17.73271209059684
$ws.Range("N5").Value = 16.62982047292599
$ws.Range("O5").Value = 16.96042918365005
$ws.Range("C6").Value = 4.587811108341602
$ws.Range("D6").Value = 2.686071142315282
$ws.Range("E6").Value = 9.903033097455941
$ws.Range("F6").Value = 19.46952953201829
$ws.Range("G6").Value = 22.3345546493093
$ws.Range("H6").Value = 11.61300691742615
$ws.Range("M6").Value = 17.69810900141879
$ws.Range("N6").Value = 16.61735767763066
$ws.Range("O6").Value = 16.95179234847568
$ws.Range("C7").Value = 4.63985534976544
$ws.Range("D7").Value = 2.687737130977207
$ws.Range("E7").Value = 9.847148966460157
$ws.Range("F7").Value = 19.58111259119075
$ws.Range("G7").Value = 22.59593138995519
$ws.Range("H7").Value = 11.62237187348921
$ws.Range("M7").Value = 17.93695881323302
$ws.Range("N7").Value = 16.70384869011264
$ws.Range("O7").Value = 17.01218299487144
$ws.Range("C8").Value = 4.865975381994392
$ws.Range("D8").Value = 2.694993414950076
$ws.Range("E8").Value = 9.611994682128888
$ws.Range("F8").Value = 20.08186387968247
$ws.Range("G8").Value = 23.73400874614285
$ws.Range("H8").Value = 11.67055433466331
$ws.Range("M8").Value = 18.95594771559181
$ws.Range("N8").Value = 17.08449682360818
$ws.Range("O8").Value = 17.28914075010145
$ws.Range("C9").Value = 5.291731728427412
$ws.Range("D9").Value = 2.712014578552664
$ws.Range("E9").Value = 9.191710091668625
$ws.Range("F9").Value = 21.08375229268482
$ws.Range("G9").Value = 25.88466301249707
$ws.Range("H9").Value = 11.7892536302659
$ws.Range("M9").Value = 20.81203555523141
$ws.Range("N9").Value = 17.8211873955653
$ws.Range("O9").Value = 17.86529558418908
$ws.Range("C10").Value = 5.589987862462325
$ws.Range("D10").Value = 2.726743433156718
$ws.Range("E10").Value = 8.907543608317127
$ws.Range("F10").Value = 21.82312866677711
$ws.Range("G10").Value = 27.39542793039741
$ws.Range("H10").Value = 11.89075318945811
$ws.Range("M10").Value = 22.07728721949012
$ws.Range("N10").Value = 18.351285742766
$ws.Range("O10").Value = 18.30446592041041
$ws.Range("C11").Value = 5.721849196745278
$ws.Range("D11").Value = 2.733545955857736
$ws.Range("E11").Value = 8.78355321051586
$ws.Range("F11").Value = 22.15841634927532
$ws.Range("G11").Value = 28.0641793941941
$ws.Range("H11").Value = 11.93991148725811
$ws.Range("M11").Value = 22.62956610068318
$ws.Range("N11").Value = 18.58898882254293
$ws.Range("O11").Value = 18.50677627309666
$ws.Range("C12").Value = 5.771185884217644
$ws.Range("D12").Value = 2.736135667477355
$ws.Range("E12").Value = 8.737355782107988
$ws.Range("F12").Value = 22.28508816263104
$ws.Range("G12").Value = 28.31451248338952
$ws.Range("H12").Value = 11.9589443359829
$ws.Range("M12").Value = 22.8352189349635
$ws.Range("N12").Value = 18.67842529277015
$ws.Range("O12").Value = 18.58367251751546
$ws.Range("C13").Value = 5.760587562954515
$ws.Range("D13").Value = 2.735577330994862
$ws.Range("E13").Value = 8.747271706627581
$ws.Range("F13").Value = 22.25782239921756
$ws.Range("G13").Value = 28.26073175402013
$ws.Range("H13").Value = 11.95482688971169
$ws.Range("M13").Value = 22.79108475598149
$ws.Range("N13").Value = 18.65919044332294
$ws.Range("O13").Value = 18.56710006036413
$ws.Range("C14").Value = 5.72592035867021
$ws.Range("D14").Value = 2.733758743118356
$ws.Range("E14").Value = 8.779737410598026
$ws.Range("F14").Value = 22.16884440352792
$ws.Range("G14").Value = 28.08483405078118
$ws.Range("H14").Value = 11.94146903012463
$ws.Range("M14").Value = 22.64655578821484
$ws.Range("N14").Value = 18.59635882136255
$ws.Range("O14").Value = 18.51309729090773
$ws.Range("C15").Value = 5.704606723881892
$ws.Range("D15").Value = 2.732646566520032
$ws.Range("E15").Value = 8.79972180392711
$ws.Range("F15").Value = 22.1143003818862
$ws.Range("G15").Value = 27.97670601616608
$ws.Range("H15").Value = 11.93334098835371
$ws.Range("M15").Value = 22.55757037240797
$ws.Range("N15").Value = 18.5577952091048
$ws.Range("O15").Value = 18.48005394302955
$ws.Range("C16").Value = 5.581289417246986
$ws.Range("D16").Value = 2.726300848443845
$ws.Range("E16").Value = 8.915752544878064
$ws.Range("F16").Value = 21.80118321636273
$ws.Range("G16").Value = 27.35132922224237
$ws.Range("H16").Value = 11.88759967927127
$ws.Range("M16").Value = 22.04071390737825
$ws.Range("N16").Value = 18.3356749580748
$ws.Range("O16").Value = 18.29128867231091
$ws.Range("C17").Value = 5.504625003854366
$ws.Range("D17").Value = 2.72243347269696
$ws.Range("E17").Value = 8.988282708783837
$ws.Range("F17").Value = 21.60872152477603
$ws.Range("G17").Value = 26.96275542289361
$ws.Range("H17").Value = 11.86029552273561
$ws.Range("M17").Value = 21.717572028588
$ws.Range("N17").Value = 18.19847048428272
$ws.Range("O17").Value = 18.17607938050295
$ws.Range("C18").Value = 5.460172850194825
$ws.Range("D18").Value = 2.72021874743915
$ws.Range("E18").Value = 9.030497198130409
$ws.Range("F18").Value = 21.49793111100154
$ws.Range("G18").Value = 26.73752920756498
$ws.Range("H18").Value = 11.84487266797238
$ws.Range("M18").Value = 21.5295237878466
$ws.Range("N18").Value = 18.11923480024647
$ws.Range("O18").Value = 18.11005676996613
$ws.Range("C19").Value = 5.44506237471732
$ws.Range("D19").Value = 2.719470572463591
$ws.Range("E19").Value = 9.04487579146204
$ws.Range("F19").Value = 21.46040810646109
$ws.Range("G19").Value = 26.66098261471905
$ws.Range("H19").Value = 11.83969949138133
$ws.Range("M19").Value = 21.4654829302548
$ws.Range("N19").Value = 18.09235467907817
$ws.Range("O19").Value = 18.08774672935043
$ws.Range("C20").Value = 5.512823352488273
$ws.Range("D20").Value = 2.722844165848057
$ws.Range("E20").Value = 8.980510334356236
$ws.Range("F20").Value = 21.62921992684738
$ws.Range("G20").Value = 27.00430058978196
$ws.Range("H20").Value = 11.86317301614703
$ws.Range("M20").Value = 21.75219826773694
$ws.Range("N20").Value = 18.21310978655373
$ws.Range("O20").Value = 18.188319053199
$ws.Range("C21").Value = 5.73611949114679
$ws.Range("D21").Value = 2.734292541215245
$ws.Range("E21").Value = 8.770180986010107
$ws.Range("F21").Value = 22.19498848358244
$ws.Range("G21").Value = 28.13658023092776
$ws.Range("H21").Value = 11.94538131878677
$ws.Range("M21").Value = 22.68910295144796
$ws.Range("N21").Value = 18.61483026758165
$ws.Range("O21").Value = 18.52895207637706
$ws.Range("C22").Value = 5.878562463133065
$ws.Range("D22").Value = 2.741854369740094
$ws.Range("E22").Value = 8.637117781299558
$ws.Range("F22").Value = 22.56297449296282
$ws.Range("G22").Value = 28.85954322848511
$ws.Range("H22").Value = 12.00153772833907
$ws.Range("M22").Value = 23.28107307409887
$ws.Range("N22").Value = 18.87398264194251
$ws.Range("O22").Value = 18.75320537687561
$ws.Range("C23").Value = 5.802871881720349
$ws.Range("D23").Value = 2.737811513462056
$ws.Range("E23").Value = 8.707734887623523
$ws.Range("F23").Value = 22.36678105197642
$ws.Range("G23").Value = 28.47531814565714
$ws.Range("H23").Value = 11.97134784083907
$ws.Range("M23").Value = 22.96702883954219
$ws.Range("N23").Value = 18.73600463685318
$ws.Range("O23").Value = 18.63339352916113
$ws.Range("C24").Value = 5.509118048978546
$ws.Range("D24").Value = 2.722658464388652
$ws.Range("E24").Value = 8.984022617211078
$ws.Range("F24").Value = 21.61995302938958
$ws.Range("G24").Value = 26.9855237005764
$ws.Range("H24").Value = 11.86187124485185
$ws.Range("M24").Value = 21.73655080413191
$ws.Range("N24").Value = 18.20649245748161
$ws.Range("O24").Value = 18.18278482843748
$ws.Range("C25").Value = 5.178861000523603
$ws.Range("D25").Value = 2.706687516404662
$ws.Range("E25").Value = 9.301064664724443
$ws.Range("F25").Value = 20.81151187197308
$ws.Range("G25").Value = 25.31373075300049
$ws.Range("H25").Value = 11.75459482210793
$ws.Range("M25").Value = 20.32645524317606
$ws.Range("N25").Value = 17.62348762088206
$ws.Range("O25").Value = 17.70632989658833
